$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style of the existing header row (bold/centered/bordered)
# by copying the formatting from the adjacent header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's W/L/T record for every player row (2 through 47).
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 103  # AD
    $ws.Cells.Item($r, 31).Value = 59   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
